# ejercicio2_normalizacion.xlsx — "Terminado de corregir 2"
#
# On sheet "3FN" (the 3rd-normal-form worksheet), the "pot." (power) column
# is redundant against "modelo" (it only depends on the model, not on the
# full matricula+fecha key), so it is split out into its own small
# lookup table (modelo -> pot.) appended below the existing
# (precio -> ciudad) lookup table, mirroring that table's layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3FN")
[void]$ws.Activate()

# ---------------------------------------------------------------
# 1) Header row 9: drop the "pot." column; "fecha"/"precio" shift left
# ---------------------------------------------------------------
$ws.Range("C9").Value2 = "fecha"
$ws.Range("D9").Value2 = "precio"
$ws.Range("E9").Clear()

# ---------------------------------------------------------------
# 2) Data rows 10-15: drop the power values (C), shift fecha (was D)
#    into C and precio (was E) into D, taking their formatting along.
# ---------------------------------------------------------------
[void]$ws.Range("D10:D15").Copy()
[void]$ws.Range("C10:C15").PasteSpecial(-4122)   # xlPasteFormats
[void]$ws.Range("E10:E15").Copy()
[void]$ws.Range("D10:D15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C10").Value2 = 36972
$ws.Range("D10").Value2 = 15000000

$ws.Range("C11").Value2 = 37025
$ws.Range("D11").Value2 = 6000000

$ws.Range("C12").Value2 = 36845
$ws.Range("D12").Value2 = 14000000

$ws.Range("C13").Value2 = 34856
$ws.Range("D13").Value2 = 2500000

$ws.Range("C14").Value2 = 32408
$ws.Range("D14").Value2 = 1800000

$ws.Range("C15").Value2 = 32408
$ws.Range("D15").Value2 = 1800000

[void]$ws.Range("E10:E15").Clear()

# ---------------------------------------------------------------
# 3) New lookup table at rows 32-38: modelo -> pot.
#    (same visual style as the existing precio -> ciudad table
#    at rows 25-30: bold header row, bordered data rows)
# ---------------------------------------------------------------
[void]$ws.Range("A25:B25").Copy()
[void]$ws.Range("A32:B32").PasteSpecial(-4122)   # xlPasteFormats
[void]$ws.Range("A26:B30").Copy()
[void]$ws.Range("A33:B38").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A32").Value2 = "modelo"
$ws.Range("B32").Value2 = "pot."

$ws.Range("A33").Value2 = "Senator Luxe Top"
$ws.Range("B33").Value2 = 125

$ws.Range("A34").Value2 = "Espace VX"
$ws.Range("B34").Value2 = 90

$ws.Range("A35").Value2 = "Senator Luxe Top"
$ws.Range("B35").Value2 = 125

$ws.Range("A36").Value2 = "Xara JR"
$ws.Range("B36").Value2 = 65

$ws.Range("A37").Value2 = "Fiesta 1000"
$ws.Range("B37").Value2 = 50

$ws.Range("A38").Value2 = "Fiesta 1000"
$ws.Range("B38").Value2 = 50

# ---------------------------------------------------------------
# 4) View state: scroll so row 4 is at the top, with F10 selected
# ---------------------------------------------------------------
[void]$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F10").Select()

# ---------------------------------------------------------------
# 5) Page setup: portrait, paper size 9 (A4)
# ---------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
